# Applies the "Attendance reading from excel file completed / Config changed"
# commit to Data/Config.xlsx.
#
# Summary of the content-level changes (ignoring pure cosmetic/XML noise that
# Excel re-serialises on every save):
#
#   Settings sheet:
#     - Old rows 2 & 4 (OrchestratorQueueName / logF_BusinessProcessName)
#       removed, replaced by a new 6-row block of attendance-framework
#       settings (rows 2-7).
#
#   Constants sheet:
#     - MaxRetryNumber (B3) changed from 0 to 2.
#     - Six new rows (25-30) appended describing the attendance column
#       layout.
#
#   Assets sheet:
#     - No data change (header only).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")

# Wipe out the old rows 2-4 (A:C) completely - the new data only uses
# columns A and B, so column C must end up blank for these rows.
$settings.Range("A2:C7").ClearContents()

$settings.Range("A2").Value = "Attendance Format"
$settings.Range("B2").Value = "1/0"

$settings.Range("A3").Value = "Attendance Threshold"
$settings.Range("B3").Value = 75

$settings.Range("A4").Value = "Input Directory"
$settings.Range("B4").Value = "C:\Users\Samyak\Documents\UiPath\UiPath-Attendance-Framework-UAF\Data\Input"

$settings.Range("A5").Value = "Output Directory"
$settings.Range("B5").Value = "C:\Users\Samyak\Documents\UiPath\UiPath-Attendance-Framework-UAF\Data\Output"

$settings.Range("A6").Value = "Number_Batches"
$settings.Range("B6").Value = 3

$settings.Range("A7").Value = "Organization"
$settings.Range("B7").Value = "SUAS"

$settings.Range("B3").Select()

# ---------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")

# MaxRetryNumber value changed 0 -> 2
$constants.Range("B3").Value = 2

# New rows describing the attendance column configuration
$constants.Range("A25").Value = "Number Sessions Cell"
$constants.Range("B25").Value = "B1"

$constants.Range("A26").Value = "Attendance Column Fields"
$constants.Range("B26").Value = "A:C"
$constants.Range("C26").Value = "Column Names to fetch from attendance file"

$constants.Range("A27").Value = "Attendance Column Field Count"
$constants.Range("B27").Value = 3

$constants.Range("A28").Value = "Headers Row"
$constants.Range("B28").Value = 4
$constants.Range("C28").Value = "Row number containing column names"

$constants.Range("A29").Value = "Start Date Column"
$constants.Range("B29").Value = "D"

$constants.Range("A30").Value = "Subject Name Cell"
$constants.Range("B30").Value = "E1"

# Match the vertical-top alignment style applied to the new block in the
# authored workbook (rows 25-29; row 30 stays in the default/normal style).
# Only the cells that actually carry data get the style (mirrors the source
# file, which leaves B28 in the default style while its neighbours A28/C28
# are styled).
$constants.Range("A25:B25").VerticalAlignment = -4160
$constants.Range("A26:C26").VerticalAlignment = -4160
$constants.Range("A27:C27").VerticalAlignment = -4160
$constants.Range("A28").VerticalAlignment = -4160
$constants.Range("C28").VerticalAlignment = -4160
$constants.Range("A29:B29").VerticalAlignment = -4160

$constants.Range("A4").Select()
$constants.Range("B28").Select()

# ---------------------------------------------------------------------
# Assets sheet - header text only, no underlying data changed (the shared
# string table was simply reshuffled upstream).
# ---------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")
$assets.Range("B1").Value = "Asset"
$assets.Range("C1").Value = "Description (Assets will always overwrite other config)"
$assets.Range("A2").Select()

$settings.Activate()
